$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample data rows (2,3,5,6,7) entirely, leaving only the header
# row (1) and row 4 (whose content is cleared too, keeping just the style
# on I4) behind, matching the "working exercise" template state.
$ws.Rows.Item(2).ClearContents()
$ws.Rows.Item(3).ClearContents()
$ws.Rows.Item(5).ClearContents()
$ws.Rows.Item(6).ClearContents()
$ws.Rows.Item(7).ClearContents()

# Row 4 keeps its row, but every value is cleared - including I4, whose
# highlight style (s="1") is preserved for the blank cell.
$ws.Range("A4:H4").ClearContents()
$ws.Range("I4").ClearContents()

# Reflect the updated selection left behind in the sheet view.
$ws.Range("A2:I7").Select()
